# Apply updated crypto price/volume figures (diff-derived).
# Leading apostrophe forces text storage so numeric-looking strings
# (e.g. "1.001", "30.344.10") are not reinterpreted as numbers/dates.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.344.10"
$ws.Range("E2").Value = "'  +0.55%  "
$ws.Range("D3").Value = "'1.933.47"
$ws.Range("E3").Value = "'  +0.49%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "'  +0.25%  "
$ws.Range("D5").Value = "'250.99"
$ws.Range("E5").Value = "'  +2.17%  "
$ws.Range("D6").Value = "'0.7180"
$ws.Range("E6").Value = "'  -0.35%  "
$ws.Range("E7").Value = "'  +0.30%  "
$ws.Range("D8").Value = "'0.3271"
$ws.Range("E8").Value = "'  +0.73%  "
$ws.Range("D9").Value = "'27.53"
$ws.Range("E9").Value = "'  +4.13%  "
$ws.Range("D10").Value = "'0.07181"
$ws.Range("E10").Value = "'  +4.98%  "
$ws.Range("D11").Value = "'0.7990"
$ws.Range("E11").Value = "'  +0.48%  "
$ws.Range("D12").Value = "'0.08081"
$ws.Range("E12").Value = "'  +2.01%  "
$ws.Range("D13").Value = "'1.931.31"
$ws.Range("E13").Value = "'  +0.55%  "
$ws.Range("D14").Value = "'5.416"
$ws.Range("E14").Value = "'  +0.41%  "
$ws.Range("D15").Value = "'94.51"
$ws.Range("E15").Value = "'  +0.16%  "
$ws.Range("D16").Value = "'14.77"
$ws.Range("E16").Value = "'  +1.92%  "
$ws.Range("D17").Value = "'30.328.06"
$ws.Range("E17").Value = "'  +0.49%  "
$ws.Range("D18").Value = "'251.62"
$ws.Range("E18").Value = "'  -3.06%  "
$ws.Range("D19").Value = "'0.000008103"
$ws.Range("E19").Value = "'  +1.93%  "
$ws.Range("D20").Value = "'5.785"
$ws.Range("E20").Value = "'  -0.86%  "
$ws.Range("D21").Value = "'2.186.41"
$ws.Range("E22").Value = "'  +0.23%  "
$ws.Range("E23").Value = "'  +0.36%  "
$ws.Range("D24").Value = "'6.911"
$ws.Range("E24").Value = "'  +0.74%  "
$ws.Range("D25").Value = "'9.734"
$ws.Range("E25").Value = "'  +0.60%  "
$ws.Range("D26").Value = "'165.47"
$ws.Range("E26").Value = "'  +2.99%  "
$ws.Range("D27").Value = "'19.24"
$ws.Range("E27").Value = "'  +1.97%  "
$ws.Range("D28").Value = "'2.328"
$ws.Range("E28").Value = "'  +3.62%  "
$ws.Range("D29").Value = "'0.1285"
$ws.Range("E29").Value = "'  -3.99%  "
$ws.Range("E30").Value = "'  +0.52%  "
$ws.Range("E31").Value = "'  -0.18%  "
$ws.Range("D32").Value = "'4.421"
$ws.Range("D33").Value = "'4.204"
$ws.Range("E33").Value = "'  +0.55%  "
$ws.Range("D34").Value = "'0.05198"
$ws.Range("E34").Value = "'  +3.18%  "
$ws.Range("E35").Value = "'  +6.37%  "
$ws.Range("D36").Value = "'0.7469"
$ws.Range("E36").Value = "'  +1.46%  "
$ws.Range("D37").Value = "'2.771"
$ws.Range("E37").Value = "'  +1.46%  "
$ws.Range("D38").Value = "'0.01960"
$ws.Range("E38").Value = "'  +1.15%  "
$ws.Range("E39").Value = "'  -0.26%  "
$ws.Range("D40").Value = "'78.87"
$ws.Range("E40").Value = "'  -1.44%  "
$ws.Range("D41").Value = "'6.451"
$ws.Range("E41").Value = "'  -0.72%  "
$ws.Range("D42").Value = "'0.4528"
$ws.Range("E42").Value = "'  +2.24%  "
$ws.Range("D43").Value = "'2.024"
$ws.Range("E43").Value = "'  +1.20%  "
$ws.Range("E44").Value = "'  +0.25%  "
$ws.Range("D45").Value = "'0.8401"
$ws.Range("E45").Value = "'  +1.24%  "
$ws.Range("D46").Value = "'101.88"
$ws.Range("E46").Value = "'  -0.53%  "
$ws.Range("D47").Value = "'9.783"
$ws.Range("E47").Value = "'  +0.78%  "
$ws.Range("D48").Value = "'7.400"
$ws.Range("E48").Value = "'  +1.93%  "
$ws.Range("D49").Value = "'36.62"
$ws.Range("E49").Value = "'  +1.59%  "
$ws.Range("B50").Value = "'Cronos"
$ws.Range("C50").Value = "'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "'0.06064"
$ws.Range("E50").Value = "'  +2.78%  "
$ws.Range("B51").Value = "'Decentraland"
$ws.Range("C51").Value = "'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D51").Value = "'0.4176"
$ws.Range("E51").Value = "'  +1.92%  "
